$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated average_county_temperature (column K) with NOAA data
$ws.Range("K2").Value = 21.79166666666666
$ws.Range("K3").Value = 21.79166666666666
$ws.Range("K10").Value = -1.819444444444444
$ws.Range("K11").Value = -1.819444444444444
$ws.Range("K18").Value = 21.28240740740739
$ws.Range("K19").Value = 21.28240740740739
$ws.Range("K22").Value = 12.93898809523811
$ws.Range("K23").Value = 12.93898809523811

# Recomputed dependent COP columns (worst_whp_cop, best_whp_cop)
$ws.Range("R2").Value = 1.105721877767936
$ws.Range("S2").Value = 1.143718778908418

$ws.Range("R10").Value = 1.004851086664878
$ws.Range("S10").Value = 1.035188389617639

$ws.Range("R19").Value = 1.103333005990376
$ws.Range("S19").Value = 1.14113834478515

$ws.Range("R23").Value = 1.065614691876665
$ws.Range("S23").Value = 1.100460934966844
